# Generate Report for Handoff
#
# This updates the localization-status report for the 5th row group
# (45e4ab6d-0bbd-4b1a-a7d5-1695345c8c22 / 6116f279-879e-49d8-8bcc-a15fda7f2bfa
# file pair, table rows 8-13) on the Overview, zh-cn and de-de sheets:
#   - the "Latest HO Xliff Generate Date" / "Latest Handoff Datetime" timestamps
#     move forward a few seconds, and
#   - the "Priority" column is now set to "ht" for those rows (was blank).

$wb = $excel.ActiveWorkbook

$overview = $wb.Worksheets.Item("Overview")
$zhcn     = $wb.Worksheets.Item("zh-cn")
$dede     = $wb.Worksheets.Item("de-de")

# Latest HO Xliff Generate Date / Latest Handoff Datetime, rows 8-13.
# Overview!G8:G13 and de-de!H8:H13 share the same text; zh-cn!H8:H13 is the
# other one.
for ($r = 8; $r -le 13; $r++) {
    $overview.Range("G$r").Value = "2016-08-23 06:19:49"
    $dede.Range("H$r").Value     = "2016-08-23 06:19:49"
    $zhcn.Range("H$r").Value     = "2016-08-23 06:19:44"
}

# Priority column, rows 8-13, on zh-cn and de-de: blank -> "ht"
$zhcn.Range("E8:E13").Value = "ht"
$dede.Range("E8:E13").Value = "ht"
